# Update countries & provincias Spain
# Refresh the COVID case counters for the affected countries. A couple of
# countries (Mayotte, Etiopia) received updated figures that pushed them
# above their neighbour in the (descending, by total cases) ordering, so
# the neighbouring rows below them are re-written with the values that
# used to belong to the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Malta: plain data refresh (row 105) ---------------------------------
$ws.Range("B105").Value = 422
$ws.Range("C105").Value = 10
$ws.Range("D105").Value = 82
$ws.Range("E105").Value = 337
$ws.Range("F105").Value = 4
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 3

# --- Mayotte moves above Sri Lanka (rows 118-119) -------------------------
$ws.Range("A118").Value = "Mayotte"
$ws.Range("B118").Value = 245
$ws.Range("C118").Value = 12
$ws.Range("D118").Value = 117
$ws.Range("E118").Value = 124
$ws.Range("F118").Value = 6
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 4

$ws.Range("A119").Value = "Sri Lanka"
$ws.Range("B119").Value = 238
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 70
$ws.Range("E119").Value = 161
$ws.Range("F119").Value = 1
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 7

# --- Etiopia moves above Guayana Francesa / Gabon / Aruba / Tanzania /
#     Monaco (rows 136-141) -------------------------------------------------
$ws.Range("A136").Value = "Etiopia"
$ws.Range("B136").Value = 96
$ws.Range("C136").Value = 4
$ws.Range("D136").Value = 15
$ws.Range("E136").Value = 78
$ws.Range("F136").Value = 1
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 3

$ws.Range("A137").Value = "Guayana Francesa"
$ws.Range("B137").Value = 96
$ws.Range("C137").Value = 10
$ws.Range("D137").Value = 61
$ws.Range("E137").Value = 35
$ws.Range("F137").Value = 2
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

$ws.Range("A138").Value = "Gabon"
$ws.Range("B138").Value = 95
$ws.Range("C138").Value = 15
$ws.Range("D138").Value = 6
$ws.Range("E138").Value = 88
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 1

$ws.Range("A139").Value = "Aruba"
$ws.Range("B139").Value = 95
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 39
$ws.Range("E139").Value = 54
$ws.Range("F139").Value = 1
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 2

$ws.Range("A140").Value = "Tanzania"
$ws.Range("B140").Value = 94
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 11
$ws.Range("E140").Value = 79
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 4

$ws.Range("A141").Value = "Monaco"
$ws.Range("B141").Value = 93
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 12
$ws.Range("E141").Value = 78
$ws.Range("F141").Value = 2
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 3
